# Applies the "categorias" service rows + view/cosmetic tweaks described by
# the diff to 01DiseñoBDD.xlsx

$wb = $excel.ActiveWorkbook

$wsTablas = $wb.Worksheets.Item("TABLAS BDD")
$wsServicios = $wb.Worksheets.Item("SERVICIOS")

# --- New service rows on SERVICIOS sheet (categorias CRUD endpoints) -------
$wsServicios.Range("B13").Value = "POST"
$wsServicios.Range("C13").Value = "/categorias/crear"
$wsServicios.Range("D13").Value = "crea o inserta una nueva categoria en la tabla categoria"

$wsServicios.Range("B14").Value = "PUT"
$wsServicios.Range("C14").Value = "/categorias/actualizar"
$wsServicios.Range("D14").Value = "actualiza o modifica una categoria"

$wsServicios.Range("B15").Value = "GET"
$wsServicios.Range("C15").Value = "/categorias/recuperar"
$wsServicios.Range("D15").Value = "recupera todas las categorias registradas en la tabla."

# C14 picks up the same style as the rest of the column (loses the stray
# underline formatting it had before data was entered).
$wsServicios.Range("C14").Font.Underline = $false

# --- Column width tweak on SERVICIOS sheet ----------------------------------
$wsServicios.Columns.Item(1).ColumnWidth = 3.5

# --- Selection / scroll position bookkeeping --------------------------------
# Scroll "TABLAS BDD" so row 37 becomes the top-left visible row, without
# disturbing its existing cell selection (G61).
$wsTablas.Activate()
$excel.Goto($wsTablas.Range("A37"), $true) | Out-Null
$wsTablas.Range("G61").Select() | Out-Null

# SERVICIOS stays the active sheet/tab, with the selection moved to D18.
$wsServicios.Activate()
$wsServicios.Range("D18").Select() | Out-Null
